# Add a "Save" column (H) to the s_vals sheet, matching the header style
# used by the other header cells and filling in the per-row save values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy formatting from an existing header cell (G1) so it
# gets the same bold font / border / center-top alignment style, then set
# its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Values for H2:H12 ("Save" flag per row), plain numbers with no special
# style, matching column G's un-styled data cells.
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
